$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value2 = "14/03/2023"

$data = @(
  @("B2", 0),
  @("C2", 522),
  @("D2", 523),
  @("E2", 0),
  @("F2", 1),
  @("G2", 11),
  @("H2", 0),
  @("I2", 525.2),
  @("J2", -0.4188880426504227),
  @("B3", 0),
  @("C3", 212),
  @("D3", 214),
  @("E3", 0),
  @("F3", 2),
  @("G3", 8),
  @("H3", 0),
  @("I3", 236),
  @("J3", -9.322033898305083),
  @("B4", 0),
  @("C4", 6),
  @("D4", 7),
  @("E4", 1),
  @("F4", 0),
  @("G4", 0),
  @("H4", 0),
  @("I4", 6),
  @("J4", 16.66666666666667),
  @("B5", 0),
  @("C5", 266),
  @("D5", 271),
  @("E5", 2),
  @("F5", 3),
  @("G5", 0),
  @("H5", 0),
  @("I5", 96),
  @("J5", 182.2916666666667),
  @("B6", 0),
  @("C6", 28),
  @("D6", 28),
  @("E6", 0),
  @("F6", 0),
  @("G6", 0),
  @("H6", 0),
  @("I6", 38),
  @("J6", -26.31578947368422),
  @("B7", 2),
  @("C7", 48),
  @("D7", 55),
  @("E7", 5),
  @("F7", 0),
  @("G7", 0),
  @("H7", 0),
  @("I7", 82),
  @("J7", -32.92682926829268),
  @("B8", 0),
  @("C8", 214),
  @("D8", 227),
  @("E8", 12),
  @("F8", 2),
  @("G8", 8),
  @("H8", 0),
  @("I8", 35),
  @("J8", 548.5714285714286),
  @("B9", 0),
  @("C9", 42),
  @("D9", 44),
  @("E9", 3),
  @("F9", 0),
  @("G9", 0),
  @("H9", 0),
  @("I9", 77),
  @("J9", -42.85714285714286),
  @("B10", 0),
  @("C10", 232),
  @("D10", 248),
  @("E10", 16),
  @("F10", 0),
  @("G10", 5),
  @("H10", 0),
  @("I10", 440),
  @("J10", -43.63636363636364),
  @("B11", 0),
  @("C11", 185),
  @("D11", 185),
  @("E11", 0),
  @("F11", 0),
  @("G11", 5),
  @("H11", 0),
  @("I11", 148),
  @("J11", 25),
  @("B12", 0),
  @("C12", 338),
  @("D12", 432),
  @("E12", 5),
  @("F12", 2),
  @("G12", 2),
  @("H12", 87),
  @("I12", 483.6),
  @("J12", -10.66997518610422),
  @("B13", 0),
  @("C13", 7),
  @("D13", 8),
  @("E13", 1),
  @("F13", 0),
  @("G13", 0),
  @("H13", 0),
  @("I13", 355),
  @("J13", -97.74647887323944),
  @("B14", 0),
  @("C14", 320),
  @("D14", 362),
  @("E14", 4),
  @("F14", 0),
  @("G14", 2),
  @("H14", 39),
  @("I14", 530),
  @("J14", -31.69811320754717),
  @("B15", 0),
  @("C15", 124),
  @("D15", 125),
  @("E15", 1),
  @("F15", 0),
  @("G15", 0),
  @("H15", 0),
  @("I15", 159),
  @("J15", -21.38364779874213),
  @("B16", 0),
  @("C16", 98),
  @("D16", 131),
  @("E16", 33),
  @("F16", 0),
  @("G16", 2),
  @("H16", 0),
  @("I16", 124),
  @("J16", 5.645161290322576),
  @("B17", 0),
  @("C17", 52),
  @("D17", 51),
  @("E17", 0),
  @("F17", 0),
  @("G17", 0),
  @("H17", 0),
  @("I17", 86),
  @("J17", -40.69767441860465),
  @("B18", 0),
  @("C18", 4),
  @("D18", 4),
  @("E18", 0),
  @("F18", 0),
  @("G18", 0),
  @("H18", 0),
  @("I18", 3),
  @("J18", 33.33333333333333),
  @("B19", 0),
  @("C19", 6),
  @("D19", 6),
  @("E19", 0),
  @("F19", 0),
  @("G19", 0),
  @("H19", 0),
  @("I19", 9),
  @("J19", -33.33333333333334),
  @("B20", 0),
  @("C20", 13),
  @("D20", 14),
  @("E20", 1),
  @("F20", 0),
  @("G20", 0),
  @("H20", 0),
  @("I20", 59),
  @("J20", -76.27118644067797)
)

foreach ($pair in $data) {
    $ws.Range($pair[0]).Value2 = $pair[1]
}
